# ARKCORR-22 Added elvis operator to drools rules.
# Update the CONDITION expressions in the "Set Due Date ..." rule rows to use
# the Elvis (safe-navigation) operator "queue?.name" instead of "queue.name",
# and fix up the Release-queue row so Rule Name / CONDITION / ACTION line up
# correctly (B31 / C31 / D31).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 27-30: Intake / Fulfill / Supervisor Approval / Executive Approval
$ws.Range("C27").Value = 'queue?.name == "Intake"'
$ws.Range("C28").Value = 'queue?.name == "Fulfill"'
$ws.Range("C29").Value = 'queue?.name == "Supervisor Approval"'
$ws.Range("C30").Value = 'queue?.name == "Executive Approval"'

# Row 31: Release queue rule
$ws.Range("B31").Value = "Set Due Date Release Queue"
$ws.Range("C31").Value = 'queue?.name == "Release"'
$ws.Range("D31").Value = "setDueDate, null"

# Match the author's on-disk view/selection state.
$ws.Activate()
$ws.Range("C31").Select()
